$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 176
$ws.Range("A3").Value = 129
$ws.Range("A4").Value = 126
$ws.Range("A5").Value = 123
$ws.Range("A6").Value = 110
$ws.Range("A7").Value = 102
$ws.Range("A8").Value = 101
$ws.Range("A9").Value = 96
$ws.Range("A10").Value = 92
$ws.Range("A11").Value = 90
